$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 268
$ws1.Range("F4").Value = 617
$ws1.Range("F5").Value = 2719
$ws1.Range("F7").Value = 207
$ws1.Range("F10").Value = 6074
$ws1.Range("F14").Value = 4951
$ws1.Range("F15").Value = 378
$ws1.Range("F16").Value = 95
$ws1.Range("F17").Value = 11
$ws1.Range("F18").Value = 2551
$ws1.Range("F19").Value = 1324
$ws1.Range("F20").Value = 1497
$ws1.Range("F21").Value = 1205
$ws1.Range("F22").Value = 281
$ws1.Range("F24").Value = 125
$ws1.Range("F25").Value = 1015
$ws1.Range("F26").Value = 223
$ws1.Range("F29").Value = 1344
$ws1.Range("F30").Value = 9
$ws1.Range("F31").Value = 2077
$ws1.Range("F32").Value = 286
$ws1.Range("F33").Value = 565
$ws1.Range("F34").Value = 61
$ws1.Range("F36").Value = 1467
$ws1.Range("F38").Value = 1011
$ws1.Range("F41").Value = 265
$ws1.Range("F42").Value = 1717
$ws1.Range("F43").Value = 2515
$ws1.Range("F45").Value = 108
$ws1.Range("F46").Value = 262
$ws1.Range("F49").Value = 90

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 11
$ws2.Range("F6").Value = 16
$ws2.Range("F7").Value = 401
$ws2.Range("F10").Value = 78
$ws2.Range("F23").Value = 332
$ws2.Range("F24").Value = 28
$ws2.Range("F37").Value = 24

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F7").Value = 565
$ws3.Range("F8").Value = 1442
$ws3.Range("F10").Value = 2408
$ws3.Range("F11").Value = 801
$ws3.Range("F12").Value = 700

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 617
$ws4.Range("F6").Value = 565
$ws4.Range("F7").Value = 2719
$ws4.Range("F8").Value = 207
$ws4.Range("F9").Value = 1442
$ws4.Range("F10").Value = 2408
$ws4.Range("F11").Value = 6075
$ws4.Range("F12").Value = 801
$ws4.Range("F13").Value = 11
$ws4.Range("F14").Value = 16
$ws4.Range("F16").Value = 4951
$ws4.Range("F17").Value = 95
$ws4.Range("F18").Value = 2551
$ws4.Range("F19").Value = 1324
$ws4.Range("F20").Value = 1497
$ws4.Range("F21").Value = 1205
$ws4.Range("F22").Value = 281
$ws4.Range("F24").Value = 125
$ws4.Range("F26").Value = 223
$ws4.Range("F28").Value = 1344
$ws4.Range("F29").Value = 2077
$ws4.Range("F30").Value = 286
$ws4.Range("F31").Value = 565
$ws4.Range("F34").Value = 1467
$ws4.Range("F40").Value = 265
$ws4.Range("F41").Value = 28
$ws4.Range("F42").Value = 1717
$ws4.Range("F43").Value = 2515
$ws4.Range("F44").Value = 108
$ws4.Range("F45").Value = 262
